$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.52 = 38886.65 pesos`n✅ 38886.65 pesos = 9.49 = 948.57 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 104.998
$ws2.Range("O10").Value = 4083.02
$ws2.Range("N12").Value = 4099.5
$ws2.Range("O12").Value = 100
